$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-09 Friday" "2026-01-10 Saturday"

Replace-Text "660×8=5280" "479×7=3353"
Replace-Text "375×3=1125" "187×7=1309"
Replace-Text "931×6=5586" "294×5=1470"
Replace-Text "913×5=4565" "754×9=6786"
Replace-Text "565×9=5085" "928×2=1856"

Replace-Text "839×2=1678" "897×7=6279"
Replace-Text "368×8=2944" "516×8=4128"
Replace-Text "588×5=2940" "977×5=4885"
Replace-Text "816×8=6528" "137×6=822"
Replace-Text "507×4=2028" "677×5=3385"

Replace-Text "869×6=5214" "152×4=608"
Replace-Text "807×3=2421" "525×6=3150"
Replace-Text "215×7=1505" "305×4=1220"
Replace-Text "832×4=3328" "687×8=5496"
Replace-Text "316×9=2844" "683×4=2732"

Replace-Text "817×6=4902" "864×6=5184"
Replace-Text "772×5=3860" "714×4=2856"
Replace-Text "549×9=4941" "907×6=5442"
Replace-Text "748×6=4488" "891×8=7128"
Replace-Text "288×4=1152" "469×2=938"

Replace-Text "341×8=2728" "875×6=5250"
Replace-Text "977×9=8793" "852×4=3408"
Replace-Text "254×5=1270" "916×7=6412"
Replace-Text "548×5=2740" "885×4=3540"
Replace-Text "679×6=4074" "731×5=3655"
